# Refresh the crypto price/volume table (GitHub Actions scrape update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that often looks numeric (thousands-dot
# separated, e.g. "70.017.18", or small decimals like "0.0000291").
# Prefix the literal with an apostrophe so Excel stores/keeps it as text
# instead of silently reparsing or rounding it as a number.

# Row 2
$ws.Range("D2").Value = "'70.017.18"
$ws.Range("E2").Value = '  +5.41%  '

# Row 3
$ws.Range("D3").Value = "'3.593.73"
$ws.Range("E3").Value = '  +4.89%  '

# Row 4
$ws.Range("E4").Value = '  +0.31%  '

# Row 5
$ws.Range("D5").Value = "'588.03"
$ws.Range("E5").Value = '  +3.38%  '

# Row 6
$ws.Range("D6").Value = "'190.96"
$ws.Range("E6").Value = '  +5.03%  '

# Row 7
$ws.Range("D7").Value = "'0.645"
$ws.Range("E7").Value = '  +1.90%  '

# Row 8
$ws.Range("D8").Value = "'3.582.83"
$ws.Range("E8").Value = '  +4.85%  '

# Row 9
$ws.Range("E9").Value = '  +0.10%  '

# Row 10
$ws.Range("D10").Value = "'0.178"
$ws.Range("E10").Value = '  -0.72%  '

# Row 11
$ws.Range("D11").Value = "'0.660"
$ws.Range("E11").Value = '  +2.51%  '

# Row 12
$ws.Range("D12").Value = "'57.86"
$ws.Range("E12").Value = '  +4.74%  '

# Row 13
$ws.Range("D13").Value = "'0.0000291"
$ws.Range("E13").Value = '  +3.68%  '

# Row 14
$ws.Range("D14").Value = "'9.72"
$ws.Range("E14").Value = '  +3.89%  '

# Row 15
$ws.Range("D15").Value = "'4.191.52"
$ws.Range("E15").Value = '  +5.50%  '

# Row 16
$ws.Range("D16").Value = "'3.623.34"
$ws.Range("E16").Value = '  +6.02%  '

# Row 17
$ws.Range("D17").Value = "'19.33"
$ws.Range("E17").Value = '  +5.01%  '

# Row 18
$ws.Range("D18").Value = "'70.302.95"
$ws.Range("E18").Value = '  +5.93%  '

# Row 19
$ws.Range("D19").Value = "'12.47"
$ws.Range("E19").Value = '  +3.68%  '

# Row 20
$ws.Range("E20").Value = '  +0.28%  '

# Row 21
$ws.Range("D21").Value = "'1.05"
$ws.Range("E21").Value = '  +3.63%  '

# Row 22
$ws.Range("D22").Value = "'495.16"
$ws.Range("E22").Value = '  +6.08%  '

# Row 23
$ws.Range("D23").Value = "'17.24"
$ws.Range("E23").Value = '  +17.97%  '

# Row 24
$ws.Range("D24").Value = "'5.37"
$ws.Range("E24").Value = '  +7.37%  '

# Row 25
$ws.Range("D25").Value = "'4.47"
$ws.Range("E25").Value = '  +7.44%  '

# Row 26
$ws.Range("D26").Value = "'90.51"
$ws.Range("E26").Value = '  +0.55%  '

# Row 27
$ws.Range("D27").Value = "'3.10"
$ws.Range("E27").Value = '  +5.20%  '

# Row 28
$ws.Range("D28").Value = "'11.11"
$ws.Range("E28").Value = '  +2.43%  '

# Row 29
$ws.Range("D29").Value = "'9.38"
$ws.Range("E29").Value = '  +5.28%  '

# Row 30
$ws.Range("D30").Value = "'32.18"
$ws.Range("E30").Value = '  +2.23%  '

# Row 31
$ws.Range("D31").Value = "'7.50"
$ws.Range("E31").Value = '  +8.38%  '

# Row 32
$ws.Range("D32").Value = "'620.60"
$ws.Range("E32").Value = '  +5.78%  '

# Row 33
$ws.Range("D33").Value = "'12.18"
$ws.Range("E33").Value = '  +4.88%  '

# Row 34
$ws.Range("D34").Value = "'0.117"
$ws.Range("E34").Value = '  +6.52%  '

# Row 35
$ws.Range("D35").Value = "'65.09"
$ws.Range("E35").Value = '  +3.95%  '

# Row 36
$ws.Range("D36").Value = "'0.0₃0814"
$ws.Range("E36").Value = '  +6.86%  '

# Row 37
$ws.Range("D37").Value = "'0.405"
$ws.Range("E37").Value = '  +5.37%  '

# Row 38
$ws.Range("B38").Value = 'Dai'
$ws.Range("C38").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = '  +0.07%  '

# Row 39
$ws.Range("B39").Value = 'InjectiveProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D39").Value = "'37.90"
$ws.Range("E39").Value = '  +3.60%  '

# Row 40
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").Value = "'0.146"
$ws.Range("E40").Value = '  +0.17%  '

# Row 41
$ws.Range("D41").Value = "'3.63"
$ws.Range("E41").Value = '  +0.69%  '

# Row 42
$ws.Range("D42").Value = "'3.312.33"
$ws.Range("E42").Value = '  +5.92%  '

# Row 43
$ws.Range("D43").Value = "'3.07"
$ws.Range("E43").Value = '  +4.57%  '

# Row 44
$ws.Range("D44").Value = "'0.0444"
$ws.Range("E44").Value = '  +4.52%  '

# Row 45
$ws.Range("D45").Value = "'2.67"
$ws.Range("E45").Value = '  +5.16%  '

# Row 46
$ws.Range("D46").Value = "'3.32"
$ws.Range("E46").Value = '  +4.01%  '

# Row 47
$ws.Range("D47").Value = "'0.137"
$ws.Range("E47").Value = '  +1.54%  '

# Row 48
$ws.Range("D48").Value = "'9.07"
$ws.Range("E48").Value = '  +5.16%  '

# Row 49
$ws.Range("D49").Value = "'2.69"
$ws.Range("E49").Value = '  -0.38%  '

# Row 50
$ws.Range("D50").Value = "'3.29"
$ws.Range("E50").Value = '  +4.36%  '

# Row 51
$ws.Range("E51").Value = '  +0.13%  '
